$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 86.66666666666667
$ws.Range("C2").Value = 13

$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 8

$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 14
